$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update Start_time / End_time ---
$info = $wb.Worksheets.Item("Info")
$info.Range("B26").Value = "Thu Nov 19 15:49:22 2020"
$info.Range("B27").Value = "Thu Nov 19 15:49:30 2020"

# --- Sheet "sessionInfo": update package versions & remove backports row ---
$si = $wb.Worksheets.Item("sessionInfo")

# Ohter_packages: here version 0.1 -> 1.0.0
$si.Range("G2").Value = "1.0.0"

# Loaded_only: magrittr version 1.5 -> 2.0.1
$si.Range("J3").Value = "2.0.1"

# Loaded_only: rprojroot version 1.3-2 -> 2.0.2
$si.Range("J10").Value = "2.0.2"

# Remove the "backports" row (was I15:J15) by shifting the remaining
# Loaded_only entries (boot / 1.3-25, previously I16:J16) up one row,
# then clearing the now-duplicate trailing row.
$si.Range("I15").Value = "boot"
$si.Range("J15").Value = "1.3-25"
$si.Range("I16:J16").ClearContents()
